{"js": "// Replace the date line and every three-digit x one-digit multiplication\n// expression in the practice sheet with the new values from the commit.\nconst replacements = [\n  [\"2024-03-06 Wednesday\", \"2024-03-07 Thursday\"],\n  [\"224\u00d74=896\", \"827\u00d75=4135\"],\n  [\"269\u00d77=1883\", \"949\u00d79=8541\"],\n  [\"934\u00d73=2802\", \"232\u00d78=1856\"],\n  [\"131\u00d74=524\", \"161\u00d73=483\"],\n  [\"565\u00d76=3390\", \"861\u00d76=5166\"],\n  [\"480\u00d74=1920\", \"397\u00d78=3176\"],\n  [\"771\u00d72=1542\", \"751\u00d76=4506\"],\n  [\"227\u00d75=1135\", \"344\u00d79=3096\"],\n  [\"358\u00d77=2506\", \"516\u00d72=1032\"],\n  [\"515\u00d74=2060\", \"859\u00d75=4295\"],\n  [\"750\u00d78=6000\", \"149\u00d75=745\"],\n  [\"249\u00d77=1743\", \"722\u00d73=2166\"],\n  [\"622\u00d79=5598\", \"391\u00d74=1564\"],\n  [\"866\u00d77=6062\", \"844\u00d72=1688\"],\n  [\"138\u00d77=966\", \"302\u00d74=1208\"],\n  [\"110\u00d73=330\", \"565\u00d74=2260\"],\n  [\"637\u00d77=4459\", \"217\u00d73=651\"],\n  [\"297\u00d74=1188\", \"891\u00d75=4455\"],\n  [\"135\u00d73=405\", \"476\u00d79=4284\"],\n  [\"483\u00d74=1932\", \"270\u00d73=810\"],\n  [\"630\u00d76=3780\", \"307\u00d76=1842\"],\n  [\"838\u00d77=5866\", \"810\u00d75=4050\"],\n  [\"986\u00d72=1972\", \"175\u00d79=1575\"],\n  [\"305\u00d76=1830\", \"581\u00d74=2324\"],\n  [\"376\u00d75=1880\", \"188\u00d75=940\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const found = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  found.load(\"items\");\n  await context.sync();\n\n  if (found.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n\n  found.items[0].insertText(newText, \"Replace\");\n  await context.sync();\n}\n", "ps1": "# Replace the date line and every three-digit x one-digit multiplication\n# expression in the practice sheet with the new values from the commit.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2024-03-06 Wednesday\", \"2024-03-07 Thursday\"),\n    @(\"224\u00d74=896\", \"827\u00d75=4135\"),\n    @(\"269\u00d77=1883\", \"949\u00d79=8541\"),\n    @(\"934\u00d73=2802\", \"232\u00d78=1856\"),\n    @(\"131\u00d74=524\", \"161\u00d73=483\"),\n    @(\"565\u00d76=3390\", \"861\u00d76=5166\"),\n    @(\"480\u00d74=1920\", \"397\u00d78=3176\"),\n    @(\"771\u00d72=1542\", \"751\u00d76=4506\"),\n    @(\"227\u00d75=1135\", \"344\u00d79=3096\"),\n    @(\"358\u00d77=2506\", \"516\u00d72=1032\"),\n    @(\"515\u00d74=2060\", \"859\u00d75=4295\"),\n    @(\"750\u00d78=6000\", \"149\u00d75=745\"),\n    @(\"249\u00d77=1743\", \"722\u00d73=2166\"),\n    @(\"622\u00d79=5598\", \"391\u00d74=1564\"),\n    @(\"866\u00d77=6062\", \"844\u00d72=1688\"),\n    @(\"138\u00d77=966\", \"302\u00d74=1208\"),\n    @(\"110\u00d73=330\", \"565\u00d74=2260\"),\n    @(\"637\u00d77=4459\", \"217\u00d73=651\"),\n    @(\"297\u00d74=1188\", \"891\u00d75=4455\"),\n    @(\"135\u00d73=405\", \"476\u00d79=4284\"),\n    @(\"483\u00d74=1932\", \"270\u00d73=810\"),\n    @(\"630\u00d76=3780\", \"307\u00d76=1842\"),\n    @(\"838\u00d77=5866\", \"810\u00d75=4050\"),\n    @(\"986\u00d72=1972\", \"175\u00d79=1575\"),\n    @(\"305\u00d76=1830\", \"581\u00d74=2324\"),\n    @(\"376\u00d75=1880\", \"188\u00d75=940\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $rng = $d.Content\n    $found = $rng.Find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n\n    if (-not $found) {\n        Write-Output \"NOT FOUND: $oldText\"\n    }\n}\n\nWrite-Output \"done\"\n"}
